# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with newly scraped values. The Price column holds plain text
# (values like "27.686.66" use dots as thousands separators, and some
# are legitimate decimals like "1.002"), so NumberFormat is forced to
# Text ("@") before each write to keep Excel from reinterpreting a
# numeric-looking string (e.g. "1.001") as a Number. The Volume(1h)
# column is already safely text (padded with spaces + a trailing "%").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.717.37'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.758.04'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.26'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4604'
$ws.Range('E7').Value = '  +7.17%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3605'
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07522'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.14'
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.100'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.81'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.030'
$ws.Range('E14').Value = '  -1.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.118'
$ws.Range('E15').Value = '  -2.97%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.747.49'
$ws.Range('E16').Value = '  -2.53%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.41'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001067'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06412'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.81'
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.831'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '27.763.24'
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.105'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.45'
$ws.Range('E26').Value = '  +3.46%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.42'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.957.92'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.095'
$ws.Range('E29').Value = '  -3.13%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '126.72'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.066'
$ws.Range('E31').Value = '  -8.66%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09240'
$ws.Range('E32').Value = '  +3.41%  '
$ws.Range('E33').Value = '  -3.56%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.564'
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '11.94'
$ws.Range('E35').Value = '  -4.78%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02303'
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2105'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06048'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6376'
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.980'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.202'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('E42').Value = '  -1.60%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.836'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.28'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5919'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.715'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '123.54'
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.959'
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.156'
$ws.Range('E49').Value = '  -3.03%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06873'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '72.54'
$ws.Range('E51').Value = '  -2.40%  '
